$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("3-9-2022 11:45:20", 60, 169153, 0, 1, 242.2610625, 4339.453375, 4339.4090625, 4338.7996875, 4340.3539375, 4338.75225, 4338.9178125, 4336.4964375, 16.05678125, 23.93606249999999, 31.96037500000001, 63.58449874999999, 127.1353625, 221.59961, 129.787545, 260.7537975, 390.68310375, 650.8980712499999, 1040.56667625, 2078.95582625, 3639.235605, 10, 6, 4, 2, 1, 0.5, 0.25, 969, 18327, 35684, 53040, 70401, 87756, 105112, 122458, 23, 38, 61, 91, 122, 243, 486, 847, 496, 996, 1493, 2487, 3976, 7944, 13906),
    @("3-9-2022 12:4:11", 61, 169153, 0, 2, 272.0392499999998, 4368.728, 4369.676937499999, 4369.562125, 4371.781, 4369.794625, 4369.190500000001, 4366.098875, 16.01034375, 23.87509375, 32.0115625, 63.32156249999998, 127.2145625, 221.55065625, 130.41603125, 261.3429062500001, 391.87328125, 652.2020625, 1043.3900625, 2084.459375, 3648.59259375, 10, 6, 4, 2, 1, 0.5, 0.25, 1088, 18563, 36042, 53520, 71007, 88486, 105963, 123427, 23, 38, 61, 91, 122, 241, 485, 844, 497, 996, 1494, 2486, 3977, 7945, 13906)
)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $rowData = $newRows[$i]
    $targetRow = 11 + $i
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($targetRow, $c + 1).Value2 = $rowData[$c]
    }
}
